$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "SecondSearchTerm" worksheet right after "SimpleSearch"
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "SecondSearchTerm"

# Populate the new sheet's data
$newSheet.Range("A2").Value = "Manhattan, NY"
$newSheet.Range("B2").Value = "The Heritage by Common"

# Match the formatting used on SimpleSearch!B2 (Menlo font style)
$ws1.Range("B2").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)

# Match column widths used on the new sheet
$newSheet.Columns.Item(1).ColumnWidth = 14.1
$newSheet.Columns.Item(2).ColumnWidth = 27.1

# Match the selected cell on the new sheet
$newSheet.Range("B2").Select() | Out-Null

# The new sheet becomes the active/selected tab
$newSheet.Activate()
